# Add a "2022-Q1" worksheet (fund-holdings detail) right before the
# "总计" (totals) sheet, and prepend a corresponding summary row to
# "总计" itself. Mirrors the other quarterly sheets' layout/number-vs-text
# typing so the new sheet matches its siblings (e.g. "2021-Q4").

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new sheet, positioned immediately before "总计".
# ---------------------------------------------------------------------
$zongji = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($zongji)
$newSheet.Name = "2022-Q1"

# Small helper to stamp the bold / centered / bordered "header-ish"
# look used on this workbook's header row and on the leading index
# column (matches column A / row 1 styling on the sibling sheets).
function Set-HeaderStyle($cell) {
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108   # xlCenter
    $cell.VerticalAlignment = -4160     # xlVAlignTop
    $cell.BorderAround(1) | Out-Null    # xlContinuous, default thin weight
}

# Helper: write a value as literal TEXT (not auto-converted to a number)
# by forcing Text number format first - matches the source data where
# numeric-looking figures ("43.33", "73.69", ...) are stored as strings.
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

# ---------------------------------------------------------------------
# 2. Header row for "2022-Q1".
# ---------------------------------------------------------------------
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$col = 2
foreach ($h in $headers) {
    $c = $newSheet.Cells.Item(1, $col)
    $c.Value = $h
    Set-HeaderStyle $c
    $col = $col + 1
}

# ---------------------------------------------------------------------
# 3. Data rows for "2022-Q1" (fund code, name, scale, total position,
#    position ratio, held market value, position rank).
# ---------------------------------------------------------------------
$rows = @(
    @("011021", "汇添富互联网核心资产六个月持有期混合A", "43.33", "73.69", "3.82", "1.6552", 5),
    @("161611", "融通内需驱动混合", "12.96", "65.58", "2.26", "0.2929", 9),
    @("001150", "融通互联网传媒灵活配置混合", "9.18", "92.66", "2.61", "0.2396", 9),
    @("004350", "汇丰晋信价值先锋股票", "4.99", "93.32", "3.39", "0.1692", 4),
    @("011022", "汇添富互联网核心资产六个月持有期混合C", "3.70", "73.69", "3.82", "0.1413", 5),
    @("013345", "富荣信息技术混合A", "1.96", "90.39", "4.97", "0.0974", 2),
    @("003655", "信达澳银新财富灵活配置混合", "11.86", "25.86", "0.75", "0.0890", 3),
    @("002291", "诺安安鑫灵活配置混合", "2.19", "81.55", "3.54", "0.0775", 7),
    @("013346", "富荣信息技术混合C", "1.44", "90.39", "4.97", "0.0716", 2),
    @("006277", "中金瑞和灵活配置混合A", "2.02", "84.11", "3.31", "0.0669", 7),
    @("011703", "中金鑫瑞优选一年持有期灵活配置混合型证券投资基金", "1.62", "72.66", "2.80", "0.0454", 7),
    @("161038", "富国新兴成长量化精选混合（LOF）", "1.13", "93.66", "1.79", "0.0202", 8),
    @("002189", "农银汇理国企改革灵活配置混合", "1.15", "54.15", "1.62", "0.0186", 10),
    @("012005", "信达澳银恒盛混合A", "1.87", "31.90", "0.82", "0.0153", 5),
    @("003717", "中银量化精选灵活配置混合A", "0.49", "90.38", "1.18", "0.0058", 7),
    @("006278", "中金瑞和灵活配置混合C", "0.16", "84.11", "3.31", "0.0053", 7),
    @("012006", "信达澳银恒盛混合C", "0.31", "31.90", "0.82", "0.0025", 5),
    @("010484", "中银量化精选灵活配置混合C", "0.01", "90.38", "1.18", "0.0001", 7)
)

$r = 2
foreach ($row in $rows) {
    $idxCell = $newSheet.Cells.Item($r, 1)
    $idxCell.Value = $r - 2
    Set-HeaderStyle $idxCell

    Set-TextValue $newSheet.Cells.Item($r, 2) $row[0]
    Set-TextValue $newSheet.Cells.Item($r, 3) $row[1]
    Set-TextValue $newSheet.Cells.Item($r, 4) $row[2]
    Set-TextValue $newSheet.Cells.Item($r, 5) $row[3]
    Set-TextValue $newSheet.Cells.Item($r, 6) $row[4]
    Set-TextValue $newSheet.Cells.Item($r, 7) $row[5]
    $newSheet.Cells.Item($r, 8).Value = $row[6]

    $r = $r + 1
}

# ---------------------------------------------------------------------
# 4. Insert the new "2022-Q1" summary row at the top of "总计"'s data
#    (row 2, right under the header) and renumber the existing index
#    column (A) so it keeps running 0,1,2,....
#
#    NOTE: re-resolve "总计" by name here rather than reusing the
#    $zongji handle captured before Worksheets.Add() - this host binds
#    worksheet variables by position, so the older handle would now
#    silently refer to the newly inserted "2022-Q1" sheet instead.
# ---------------------------------------------------------------------
$zongji = $wb.Worksheets.Item("总计")
$zongji.Rows.Item(2).Insert()

$a2 = $zongji.Cells.Item(2, 1)
$a2.Value = 0
Set-HeaderStyle $a2
$zongji.Cells.Item(2, 2).Value = "2022-Q1"
$zongji.Cells.Item(2, 3).Value = 18
$zongji.Cells.Item(2, 4).Value = 3.01

for ($row = 3; $row -le 7; $row++) {
    $zongji.Cells.Item($row, 1).Value = $row - 2
}
